$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.813.85"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "3.429.30"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'582.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").Value = "'130.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("E9").Value = "  +4.30%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "4.008.30"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "'0.0000178"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "3.425.78"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").Value = "63.780.81"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'13.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'385.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "3.565.77"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "'73.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("D27").Value = "'0.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").Value = "'7.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("D30").Value = "'7.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").Value = "'0.155"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("D33").Value = "3.455.90"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'22.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").Value = "'5.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("D37").Value = "'6.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").Value = "'164.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "'0.0774"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "'0.789"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'41.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "'1.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").Value = "'23.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.46%  "
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "2.292.87"
$ws.Range("E50").Value = "  -6.96%  "
$ws.Range("E51").Value = "  -2.17%  "
